$wb = $excel.ActiveWorkbook

# Sheet 1: '6.0-6.3'
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 58
$ws.Range("B3").Value = 71
$ws.Range("B4").Value = 83
$ws.Range("B6").Value = 103
$ws.Range("B7").Value = 112
$ws.Range("B8").Value = 120
$ws.Range("B9").Value = 127

# Sheet 2: '6.4-6.7'
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 57
$ws.Range("B3").Value = 70
$ws.Range("B4").Value = 81
$ws.Range("B6").Value = 101
$ws.Range("B7").Value = 110
$ws.Range("B8").Value = 118
$ws.Range("B9").Value = 125

# Sheet 3: '6.8-6.11'
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 55
$ws.Range("B3").Value = 68
$ws.Range("B4").Value = 80
$ws.Range("B7").Value = 108
$ws.Range("B8").Value = 116
$ws.Range("B9").Value = 123
$ws.Range("B10").Value = 130

# Sheet 4: '7.0-7.3'
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 54
$ws.Range("B3").Value = 67
$ws.Range("B4").Value = 78
$ws.Range("B5").Value = 88
$ws.Range("B7").Value = 106
$ws.Range("B8").Value = 114
$ws.Range("B9").Value = 121
$ws.Range("B10").Value = 127

# Sheet 5: '7.4-7.7'
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 53
$ws.Range("B3").Value = 66
$ws.Range("B4").Value = 77
$ws.Range("B5").Value = 87
$ws.Range("B8").Value = 112
$ws.Range("B9").Value = 119

# Sheet 6: '7.8-7.11'
$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = 52
$ws.Range("B3").Value = 64
$ws.Range("B4").Value = 75
$ws.Range("B5").Value = 85
$ws.Range("B6").Value = 94
$ws.Range("B7").Value = 102
$ws.Range("B10").Value = 124
$ws.Range("B11").Value = 130

# Sheet 7: '8.0-8.5'
$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = 51
$ws.Range("B3").Value = 63
$ws.Range("B4").Value = 74
$ws.Range("B5").Value = 84
$ws.Range("B6").Value = 92
$ws.Range("B7").Value = 100
$ws.Range("B10").Value = 121
$ws.Range("B11").Value = 127

# Sheet 8: '8.6-8.11'
$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = 50
$ws.Range("B3").Value = 61
$ws.Range("B4").Value = 72
$ws.Range("B5").Value = 82
$ws.Range("B6").Value = 90
$ws.Range("B7").Value = 98
$ws.Range("B10").Value = 119
$ws.Range("B11").Value = 125

# Sheet 9: '9.0-9.5'
$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = 48
$ws.Range("B3").Value = 60
$ws.Range("B4").Value = 70
$ws.Range("B5").Value = 80
$ws.Range("B6").Value = 89
$ws.Range("B7").Value = 96
$ws.Range("B8").Value = 104
$ws.Range("B9").Value = 110
$ws.Range("B10").Value = 117
$ws.Range("B11").Value = 123
$ws.Range("B12").Value = 128

# Sheet 10: '9.6-9.11'
$ws = $wb.Worksheets.Item(10)
$ws.Range("B2").Value = 47
$ws.Range("B3").Value = 59
$ws.Range("B4").Value = 69
$ws.Range("B5").Value = 78
$ws.Range("B6").Value = 87
$ws.Range("B7").Value = 95
$ws.Range("B8").Value = 102
$ws.Range("B9").Value = 108
$ws.Range("B10").Value = 115
$ws.Range("B11").Value = 121
$ws.Range("B12").Value = 126

# Sheet 11: '10.0-10.5'
$ws = $wb.Worksheets.Item(11)
$ws.Range("B2").Value = 46
$ws.Range("B3").Value = 57
$ws.Range("B4").Value = 68
$ws.Range("B5").Value = 77
$ws.Range("B6").Value = 85
$ws.Range("B7").Value = 93
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 107
$ws.Range("B10").Value = 113
$ws.Range("B11").Value = 119
$ws.Range("B12").Value = 124
$ws.Range("B13").Value = 129

# Sheet 12: '10.6-10.11'
$ws = $wb.Worksheets.Item(12)
$ws.Range("B3").Value = 56
$ws.Range("B4").Value = 66
$ws.Range("B5").Value = 75
$ws.Range("B6").Value = 84
$ws.Range("B7").Value = 92
$ws.Range("B8").Value = 99
$ws.Range("B9").Value = 105
$ws.Range("B10").Value = 111
$ws.Range("B11").Value = 117
$ws.Range("B12").Value = 122
$ws.Range("B13").Value = 127

# Sheet 13: '11.0-11.5'
$ws = $wb.Worksheets.Item(13)
$ws.Range("B2").Value = 44
$ws.Range("B5").Value = 74
$ws.Range("B7").Value = 90
$ws.Range("B8").Value = 97
$ws.Range("B9").Value = 104
$ws.Range("B10").Value = 110
$ws.Range("B11").Value = 115
$ws.Range("B12").Value = 121
$ws.Range("B13").Value = 126

# Sheet 14: '11.6-11.11'
$ws = $wb.Worksheets.Item(14)
$ws.Range("B2").Value = 43
$ws.Range("B3").Value = 54
$ws.Range("B6").Value = 81
$ws.Range("B7").Value = 89
$ws.Range("B8").Value = 96
$ws.Range("B9").Value = 102
$ws.Range("B10").Value = 108
$ws.Range("B11").Value = 114
$ws.Range("B12").Value = 119
$ws.Range("B13").Value = 124
$ws.Range("B14").Value = 129

# Sheet 15: '12.0-12.5'
$ws = $wb.Worksheets.Item(15)
$ws.Range("B2").Value = 42
$ws.Range("B3").Value = 53
$ws.Range("B6").Value = 80
$ws.Range("B7").Value = 88
$ws.Range("B8").Value = 95
$ws.Range("B9").Value = 101
$ws.Range("B10").Value = 107
$ws.Range("B11").Value = 113
$ws.Range("B12").Value = 118
$ws.Range("B13").Value = 123
$ws.Range("B14").Value = 128

# Sheet 16: '12.6-12.11'
$ws = $wb.Worksheets.Item(16)
$ws.Range("B2").Value = 42
$ws.Range("B3").Value = 52
$ws.Range("B4").Value = 62
$ws.Range("B6").Value = 79
$ws.Range("B7").Value = 87
$ws.Range("B8").Value = 94
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 106
$ws.Range("B11").Value = 112
$ws.Range("B12").Value = 117
$ws.Range("B13").Value = 122
$ws.Range("B14").Value = 126

# Sheet 17: '13.0-13.11'
$ws = $wb.Worksheets.Item(17)
$ws.Range("B2").Value = 41
$ws.Range("B3").Value = 51
$ws.Range("B4").Value = 61
$ws.Range("B5").Value = 70
$ws.Range("B7").Value = 85
$ws.Range("B8").Value = 92
$ws.Range("B9").Value = 99
$ws.Range("B10").Value = 105
$ws.Range("B11").Value = 110
$ws.Range("B12").Value = 115
$ws.Range("B13").Value = 120
$ws.Range("B14").Value = 125
$ws.Range("B15").Value = 129

# Sheet 18: '14.0-14.11'
$ws = $wb.Worksheets.Item(18)
$ws.Range("B3").Value = 50
$ws.Range("B4").Value = 59
$ws.Range("B5").Value = 68
$ws.Range("B6").Value = 77
$ws.Range("B8").Value = 91
$ws.Range("B9").Value = 97
$ws.Range("B10").Value = 103
$ws.Range("B11").Value = 109
$ws.Range("B12").Value = 114
$ws.Range("B13").Value = 119
$ws.Range("B14").Value = 123
$ws.Range("B15").Value = 128

# Sheet 19: '15.0-16.11'
$ws = $wb.Worksheets.Item(19)
$ws.Range("B3").Value = 48
$ws.Range("B4").Value = 58
$ws.Range("B5").Value = 67
$ws.Range("B6").Value = 75
$ws.Range("B9").Value = 96
$ws.Range("B10").Value = 102
$ws.Range("B11").Value = 107
$ws.Range("B12").Value = 113
$ws.Range("B13").Value = 117
$ws.Range("B14").Value = 122
$ws.Range("B15").Value = 126

# Sheet 20: '17.0-18.11'
$ws = $wb.Worksheets.Item(20)
$ws.Range("B3").Value = 46
$ws.Range("B4").Value = 57
$ws.Range("B5").Value = 66
$ws.Range("B6").Value = 74
$ws.Range("B7").Value = 82
$ws.Range("B11").Value = 108
$ws.Range("B13").Value = 118
$ws.Range("B15").Value = 127
